$wb = $excel.ActiveWorkbook

# Sheet ALC, row 62 (Leve Item ID 27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 87509430
$ws.Range("I62").Value = 38472820
$ws.Range("J62").Value = 178577420
$ws.Range("K62").Value = 38472820
$ws.Range("L62").Value = 178577420
$ws.Range("M62").Value = -38472196
$ws.Range("N62").Value = -178578668

# Sheet ALC, row 65 (Leve Item ID 27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 87509430
$ws.Range("I65").Value = 38472820
$ws.Range("J65").Value = 178577420
$ws.Range("K65").Value = 192364100
$ws.Range("L65").Value = 892887100
$ws.Range("M65").Value = -192360980
$ws.Range("N65").Value = -892893340

# Sheet ALC, row 103 (Leve Item ID 19909)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 62501052
$ws.Range("I103").Value = 76923990
$ws.Range("J103").Value = 1633.3334
$ws.Range("K103").Value = 230771970
$ws.Range("L103").Value = 4900.0002
$ws.Range("M103").Value = -230771384
$ws.Range("N103").Value = -6072.0002

# Sheet ALC, row 113 (Leve Item ID 27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 16668492
$ws.Range("I113").Value = 100000000
$ws.Range("J113").Value = 2190
$ws.Range("K113").Value = 100000000
$ws.Range("L113").Value = 2190
$ws.Range("M113").Value = -99996746
$ws.Range("N113").Value = -8698

# Sheet ALC, row 116 (Leve Item ID 27778)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 29177476
$ws.Range("I116").Value = 16668161
$ws.Range("J116").Value = 41686790
$ws.Range("K116").Value = 16668161
$ws.Range("L116").Value = 41686790
$ws.Range("M116").Value = -16664719
$ws.Range("N116").Value = -41693674

# Sheet ALC, row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3704922.8
$ws.Range("I132").Value = 1389.7391
$ws.Range("J132").Value = 15873675
$ws.Range("K132").Value = 4169.2173
$ws.Range("L132").Value = 47621025
$ws.Range("M132").Value = -1639.2173
$ws.Range("N132").Value = -47626085

# Sheet ALC, row 137 (Leve Item ID 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 44147664
$ws.Range("I137").Value = 12501260
$ws.Range("J137").Value = 89356820
$ws.Range("K137").Value = 37503780
$ws.Range("L137").Value = 268070460
$ws.Range("M137").Value = -37501230
$ws.Range("N137").Value = -268075560

# Sheet ALC, row 141 (Leve Item ID 44161)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3780.4285
$ws.Range("I141").Value = 3780.4285
$ws.Range("K141").Value = 11341.2855
$ws.Range("M141").Value = -6161.2855

# Sheet ARM, row 45 (Leve Item ID 27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 371686.34
$ws.Range("I45").Value = 501262.06
$ws.Range("J45").Value = 1470
$ws.Range("K45").Value = 501262.06
$ws.Range("L45").Value = 1470
$ws.Range("M45").Value = -500885.06
$ws.Range("N45").Value = -2224

# Sheet ARM, row 122 (Leve Item ID 36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1748.5405
$ws.Range("I122").Value = 1365.6154
$ws.Range("J122").Value = 2653.6365
$ws.Range("K122").Value = 4096.8462
$ws.Range("L122").Value = 7960.9095
$ws.Range("M122").Value = -1646.8462
$ws.Range("N122").Value = -12860.9095

# Sheet CRP, row 6 (Leve Item ID 2219)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1252700.6
$ws.Range("I6").Value = 2500775.8
$ws.Range("J6").Value = 4625.5
$ws.Range("K6").Value = 2500775.8
$ws.Range("L6").Value = 4625.5
$ws.Range("M6").Value = -2500662.8
$ws.Range("N6").Value = -4851.5

# Sheet CRP, row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3595168.5
$ws.Range("I31").Value = 1984875
$ws.Range("J31").Value = 7822189
$ws.Range("K31").Value = 1984875
$ws.Range("L31").Value = 7822189
$ws.Range("M31").Value = -1984580
$ws.Range("N31").Value = -7822779

# Sheet CRP, row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3595168.5
$ws.Range("I34").Value = 1984875
$ws.Range("J34").Value = 7822189
$ws.Range("K34").Value = 1984875
$ws.Range("L34").Value = 7822189
$ws.Range("M34").Value = -1984673
$ws.Range("N34").Value = -7822593

# Sheet CRP, row 74 (Leve Item ID 10636)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 29788.1
$ws.Range("J74").Value = 29788.1
$ws.Range("L74").Value = 29788.1
$ws.Range("N74").Value = -31536.1

# Sheet CRP, row 77 (Leve Item ID 10636)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 29788.1
$ws.Range("J77").Value = 29788.1
$ws.Range("L77").Value = 89364.29999999999
$ws.Range("N77").Value = -98100.29999999999

# Sheet CUL, row 5 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2872735
$ws.Range("I5").Value = 3497322.2
$ws.Range("J5").Value = 2381987.5
$ws.Range("K5").Value = 10491966.6
$ws.Range("L5").Value = 7145962.5
$ws.Range("M5").Value = -10491854.6
$ws.Range("N5").Value = -7146186.5

# Sheet CUL, row 122 (Leve Item ID 36078)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 467.4138
$ws.Range("J122").Value = 977.6667
$ws.Range("L122").Value = 8799.0003
$ws.Range("N122").Value = -13699.0003

# Sheet CUL, row 135 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2872735
$ws.Range("I135").Value = 3497322.2
$ws.Range("J135").Value = 2381987.5
$ws.Range("K135").Value = 31475899.8
$ws.Range("L135").Value = 21437887.5
$ws.Range("M135").Value = -31473364.8
$ws.Range("N135").Value = -21442957.5

# Sheet GSM, row 10 (Leve Item ID 4306)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 800
$ws.Range("I10").Value = 500
$ws.Range("J10").Value = 1100
$ws.Range("K10").Value = 500
$ws.Range("L10").Value = 1100
$ws.Range("M10").Value = -331
$ws.Range("N10").Value = -1438

# Sheet GSM, row 14 (Leve Item ID 4198)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 7607601
$ws.Range("I14").Value = 9508751
$ws.Range("J14").Value = 3000
$ws.Range("K14").Value = 9508751
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = -9508583
$ws.Range("N14").Value = -3336

# Sheet LTW, row 46 (Leve Item ID 5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 642.8182
$ws.Range("I46").Value = 755
$ws.Range("J46").Value = 446.5
$ws.Range("K46").Value = 755
$ws.Range("L46").Value = 446.5
$ws.Range("M46").Value = -567
$ws.Range("N46").Value = -822.5

# Sheet LTW, row 55 (Leve Item ID 5284)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 6410485
$ws.Range("I55").Value = 13158117
$ws.Range("J55").Value = 234.4
$ws.Range("K55").Value = 13158117
$ws.Range("L55").Value = 234.4
$ws.Range("M55").Value = -13157944
$ws.Range("N55").Value = -580.4

# Sheet LTW, row 82 (Leve Item ID 12565)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1634.1875
$ws.Range("I82").Value = 1498.25
$ws.Range("J82").Value = 1770.125
$ws.Range("K82").Value = 1498.25
$ws.Range("L82").Value = 1770.125
$ws.Range("M82").Value = -1137.25
$ws.Range("N82").Value = -2492.125

# Sheet LTW, row 85 (Leve Item ID 12565)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1634.1875
$ws.Range("I85").Value = 1498.25
$ws.Range("J85").Value = 1770.125
$ws.Range("K85").Value = 1498.25
$ws.Range("L85").Value = 1770.125
$ws.Range("M85").Value = -250.25
$ws.Range("N85").Value = -4266.125

# Sheet LTW, row 93 (Leve Item ID 19993)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 12815.046
$ws.Range("I93").Value = 4079.2222
$ws.Range("J93").Value = 18862.924
$ws.Range("K93").Value = 4079.2222
$ws.Range("L93").Value = 18862.924
$ws.Range("M93").Value = -2831.2222
$ws.Range("N93").Value = -21358.924

Write-Host "Applied all profit updates"